# Add two new header/value column pairs (I/J) mirroring the existing
# header style used by column H ("IP"), per commit "I0 and IF added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone H1's formatting (bold/centered/bordered header style) onto I1:J1
# so the new headers match the look of the existing ones exactly.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data row values
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
